# Added ifoCAST full series evaluation.
# For every data row (2-16), the nowcast (Q0 / column B) error is dropped and
# every remaining forecast-horizon error shifts one column to the left
# (Qn -> Qn-1). Rows 2-6 previously used the full B:K width, so after the
# shift a brand new Q9 (column K) value is appended. Rows 7-15 shrink by one
# column (their last populated column is cleared). Row 16 only had the
# nowcast value, so after removing it the row has no data cells left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to append in column K for rows 2-6 (the new, one-quarter-further
# out forecast horizon error that did not exist before).
$newK = @{
    2 = -0.2804276406117954
    3 = 0.1301918964218456
    4 = 0.3403798556124878
    5 = -0.01153110738878937
    6 = -0.1586151035472806
}

# Rows 2-6: full B:K rows -> shift C:K into B:J, then set K to the new value.
foreach ($row in 2..6) {
    $shifted = $ws.Range("C$row`:K$row").Value2
    $ws.Range("B$row`:J$row").Value2 = $shifted
    $ws.Range("K$row").Value2 = $newK[$row]
}

# Rows 7-15: shift the populated part of the row one column to the left and
# drop the now-unused trailing column. Last populated column index (1-based,
# A=1) shrinks from K(11) at row 7 down to C(3) at row 15.
$lastColIndex = @{
    7  = 11  # K
    8  = 10  # J
    9  = 9   # I
    10 = 8   # H
    11 = 7   # G
    12 = 6   # F
    13 = 5   # E
    14 = 4   # D
    15 = 3   # C
}

foreach ($row in 7..15) {
    $last = $lastColIndex[$row]
    $shifted = $ws.Range($ws.Cells.Item($row, 3), $ws.Cells.Item($row, $last)).Value2
    $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, $last - 1)).Value2 = $shifted
    $ws.Cells.Item($row, $last).ClearContents()
}

# Row 16 only had the nowcast (Q0) value in column B; removing it leaves the
# row with just its label in column A.
$ws.Range("B16").ClearContents()
